$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 529 (old rows 529-549 shift down
# to become 531-551), mirroring the two extra weekly records added upstream.
$ws.Rows("529:530").Insert()

# --- New row 529 ---
$ws.Range("A529").Value = 6
$ws.Range("B529").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C529").Value = "Metropolitana"
$ws.Range("D529").Value = 44939
$ws.Range("E529").Value = 13
$ws.Range("F529").Value = 100112043
$ws.Range("G529").Value = "Pepino ensalada"
$ws.Range("H529").Value = "Sin especificar"
$ws.Range("I529").Value = "Primera"
$ws.Range("J529").Value = 450
$ws.Range("K529").Value = 11000
$ws.Range("L529").Value = 12000
$ws.Range("M529").Value = 11511
$ws.Range("N529").Value = "$/caja 60 unidades"
$ws.Range("O529").Value = "Provincia de Limarí"
$ws.Range("P529").Value = 192
$ws.Range("Q529").Value = 60
$ws.Range("R529").Value = "Hortaliza"

# --- New row 530 ---
$ws.Range("A530").Value = 6
$ws.Range("B530").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C530").Value = "Metropolitana"
$ws.Range("D530").Value = 44939
$ws.Range("E530").Value = 13
$ws.Range("F530").Value = 100112043
$ws.Range("G530").Value = "Pepino ensalada"
$ws.Range("H530").Value = "Sin especificar"
$ws.Range("I530").Value = "Primera"
$ws.Range("J530").Value = 660
$ws.Range("K530").Value = 10000
$ws.Range("L530").Value = 11000
$ws.Range("M530").Value = 10333
$ws.Range("N530").Value = "$/caja 60 unidades"
$ws.Range("O530").Value = "Provincia de Quillota"
$ws.Range("P530").Value = 172
$ws.Range("Q530").Value = 60
$ws.Range("R530").Value = "Hortaliza"
